# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table with the latest scraped figures. All of these cells hold plain
# text (the source data uses "." as a thousands separator, e.g.
# "36.576.13", and the percentages keep their original padding spaces),
# so a handful of D-column values that would otherwise look like clean
# numbers to Excel ('235.85', '54.84', ...) are written with a leading
# apostrophe to force them to stay text, exactly as typing '235.85 into
# the cell would.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.576.13'
$ws.Range('E2').Value = '  -2.15%  '
$ws.Range('D3').Value = '1.995.80'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''235.85'
$ws.Range('E5').Value = '  -9.25%  '
$ws.Range('E6').Value = '  -2.67%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '''54.84'
$ws.Range('E8').Value = '  -2.76%  '
$ws.Range('D9').Value = '''0.371'
$ws.Range('E9').Value = '  -4.20%  '
$ws.Range('D10').Value = '''58.27'
$ws.Range('E10').Value = '  +2.77%  '
$ws.Range('E11').Value = '  -3.47%  '
$ws.Range('E12').Value = '  -3.05%  '
$ws.Range('D13').Value = '''14.23'
$ws.Range('E13').Value = '  -0.55%  '
$ws.Range('D14').Value = '2.288.58'
$ws.Range('E14').Value = '  -1.07%  '
$ws.Range('D15').Value = '''20.42'
$ws.Range('E15').Value = '  -2.91%  '
$ws.Range('E16').Value = '  -5.95%  '
$ws.Range('E17').Value = '  -3.70%  '
$ws.Range('D18').Value = '2.003.85'
$ws.Range('E18').Value = '  -0.81%  '
$ws.Range('D19').Value = '36.510.73'
$ws.Range('E19').Value = '  -2.20%  '
$ws.Range('D20').Value = '''67.86'
$ws.Range('E20').Value = '  -2.95%  '
$ws.Range('E21').Value = '  -4.10%  '
$ws.Range('D22').Value = '''5.28'
$ws.Range('E22').Value = '  +1.71%  '
$ws.Range('D23').Value = '''222.06'
$ws.Range('E23').Value = '  -3.10%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('D26').Value = '''2.40'
$ws.Range('E26').Value = '  -9.46%  '
$ws.Range('D27').Value = '''162.28'
$ws.Range('E27').Value = '  -1.47%  '
$ws.Range('D28').Value = '''8.67'
$ws.Range('E28').Value = '  -3.77%  '
$ws.Range('E29').Value = '  -2.97%  '
$ws.Range('D30').Value = '''18.89'
$ws.Range('E30').Value = '  -5.84%  '
$ws.Range('E31').Value = '  +0.68%  '
$ws.Range('E32').Value = '  -3.19%  '
$ws.Range('E33').Value = '  -6.39%  '
$ws.Range('E34').Value = '  -6.54%  '
$ws.Range('D35').Value = '''4.26'
$ws.Range('E35').Value = '  -7.14%  '
$ws.Range('D36').Value = '''2.34'
$ws.Range('E36').Value = '  -2.49%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').Value = '''3.35'
$ws.Range('E38').Value = '  -0.87%  '
$ws.Range('E39').Value = '  -2.94%  '
$ws.Range('D40').Value = '''5.64'
$ws.Range('E40').Value = '  +7.32%  '
$ws.Range('E41').Value = '  -1.63%  '
$ws.Range('D42').Value = '''0.0946'
$ws.Range('E42').Value = '  +1.10%  '
$ws.Range('D43').Value = '1.452.70'
$ws.Range('E43').Value = '  +3.21%  '
$ws.Range('E44').Value = '  -5.45%  '
$ws.Range('E45').Value = '  -8.34%  '
$ws.Range('D46').Value = '''89.16'
$ws.Range('E46').Value = '  -1.41%  '
$ws.Range('E47').Value = '  -3.66%  '
$ws.Range('E48').Value = '  -3.20%  '
$ws.Range('E49').Value = '  -1.01%  '
$ws.Range('E50').Value = '  -3.84%  '
$ws.Range('D51').Value = '''3.73'
$ws.Range('E51').Value = '  +8.13%  '
